# Update the "想去人数" (want-to-go count) values in column F for the
# sheets that list individual events ("展览" and "全部类型"). Both sheets
# carry the same rows, so the same F-column updates apply to each.

$wb = $excel.ActiveWorkbook

$updates = @{
    2  = 204
    3  = 436
    4  = 12823
    5  = 1332
    6  = 188
    9  = 165
    11 = 463
    12 = 63
    16 = 392
    17 = 5471
    19 = 35
    20 = 954
    21 = 27
    23 = 111
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Range("F$row").Value = $updates[$row]
    }
}
